# Compare: accept a compare method — VQSR filter column (E) gets real
# values instead of the old placeholder a/b strings, and a new score/QD
# column (F) is populated alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: filter name + score
$ws.Range("E2").Value = "VQSR"
$ws.Range("F2").Value = 0.3

# Row 3: filter result (new)
$ws.Range("E3").Value = "PASS"

# Row 4: filter result + score
$ws.Range("E4").Value = "PASS"
$ws.Range("F4").Value = 0.004

# Row 5: filter result + score
$ws.Range("E5").Value = "PASS"
$ws.Range("F5").Value = 0.0003

# Move the live selection to the newly-populated E4:F4 block
$ws.Range("E4:F4").Select()
